$wb = $excel.ActiveWorkbook

$cypherGiantSchnauzer = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN['Giant Schnauzer']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study"

$wsStat = $wb.Worksheets.Item("StatOutput")
$wsStat.Range("A2").NumberFormat = "@"
$wsStat.Range("A2").Value = "19"
$wsStat.Range("A2").ClearFormats()

$wsStatMsg = $wb.Worksheets.Item("StatOutput_Message")
$wsStatMsg.Range("A18").Value = $cypherGiantSchnauzer
